# The Python export only writes a table/field row to the workbook when the
# corresponding checkbox is checked. Table "T2" and fields "T1F1", "T2F1",
# "T2F2" are no longer checked, so their rows must be cleared. The one
# remaining checked field ("T1F2") moves up to row 1 of Field_Names.

$wb = $excel.ActiveWorkbook

$wsTables = $wb.Worksheets.Item("Table_Names")
$wsTables.Range("A2").ClearContents()

$wsFields = $wb.Worksheets.Item("Field_Names")
$wsFields.Range("A1").Value = "T1F2"
$wsFields.Range("A2").ClearContents()
$wsFields.Range("A3").ClearContents()
$wsFields.Range("A4").ClearContents()
